$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 1.39
$ws.Range("D4").Value = 1.33
$ws.Range("G4").Value = 1.03
$ws.Range("D7").Value = 1.67
